# edit.ps1 - apply "added a bit about acclimation" commit
#
# Summary of the change (per the target diff):
#   - Insert a new slide at position 2: "Let's talk about acclimation"
#     (Section Header layout: title + empty body placeholder)
#   - Insert a new slide at position 3: "Seems pretty cut and dry!"
#     (Title Only layout: title + a cropped copy of the acclimation-types
#     figure + a small "Atkin & Tjoelker (2003)" citation textbox)
#   - On the (now 4th) slide "It turns out acclimation can be quite
#     variable" add a small "Yamori et al. (2014)" citation textbox
#
# NOTE: all Shape position/size properties (Left/Top/Width/Height) and the
# arguments to AddTextbox/AddPicture/etc. are expressed in POINTS, while the
# figures transcribed from the OOXML diff are in EMU. 1 pt = 12700 EMU.

$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) New slide at position 2 - "Let's talk about acclimation"
#    Layout 33 = ppLayoutSectionHeader ("Section Header"), matching the
#    layout used by the other title-only text slides in this deck.
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Add(2, 33)
$slide2.Shapes.Placeholders.Item(1).TextFrame.TextRange.Text = "Let's talk about acclimation"
# The body placeholder (idx 2) is left empty, same as the target slide.

# ---------------------------------------------------------------------
# 2) New slide at position 3 - "Seems pretty cut and dry!"
#    Layout 11 = ppLayoutTitleOnly ("Title Only").
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Add(3, 11)
$slide3.Shapes.Placeholders.Item(1).TextFrame.TextRange.Text = "Seems pretty cut and dry!"

# Picture: reuse the acclimation-types figure already embedded in the
# deck (the picture on the "Figures for ease of viewing" slide that used
# to be slide 8, now pushed down to slide 10) via Copy/Paste so the
# existing image data is referenced rather than duplicated.
$figSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $cand = $p.Slides.Item($i)
    if ($cand.Shapes.Count -eq 1 -and $cand.Shapes.Item(1).Type -eq 13) {
        $w = $cand.Shapes.Item(1).Width
        if ([Math]::Round($w) -eq 265) {
            $figSlide = $cand
            break
        }
    }
}
$figSlide.Shapes.Item(1).Copy()
$pastedPics = $slide3.Shapes.Paste()
$pic = $pastedPics.Item(1)
$pic.PictureFormat.CropBottom = 193.9045725
$pic.Left = 4232987 / $EMU_PER_PT
$pic.Top = 1690688 / $EMU_PER_PT
$pic.Width = 3371183 / $EMU_PER_PT
$pic.Height = 4558352 / $EMU_PER_PT

# Citation textbox: "Atkin & Tjoelker (2003)"
$tb1 = $slide3.Shapes.AddTextbox(1, 9850140 / $EMU_PER_PT, 6487033 / $EMU_PER_PT, 2341860 / $EMU_PER_PT, 369332 / $EMU_PER_PT)
$tb1.TextFrame.TextRange.Text = "Atkin & "
$r1b = $tb1.TextFrame.TextRange.InsertAfter("Tjoelker")
$r1c = $r1b.InsertAfter(" (2003)")

# ---------------------------------------------------------------------
# 3) "It turns out acclimation can be quite variable" slide (now pushed
#    down to position 4) gains a "Yamori et al. (2014)" citation textbox.
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$tb2 = $slide4.Shapes.AddTextbox(1, 10180231 / $EMU_PER_PT, 6488668 / $EMU_PER_PT, 2011769 / $EMU_PER_PT, 369332 / $EMU_PER_PT)
$tb2.TextFrame.TextRange.Text = "Yamori"
$r2b = $tb2.TextFrame.TextRange.InsertAfter(" et al. (2014)")

Write-Output "done"
